$d = $word.ActiveDocument
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$rng = $ftr.Range
$tbl = $d.Tables.Add($rng, 1, 3)
Write-Output "before count=$($ftr.Range.Paragraphs.Count)"
$tbl.Range.InsertParagraphAfter()
Write-Output "after count=$($ftr.Range.Paragraphs.Count)"
for ($i=1; $i -le $ftr.Range.Paragraphs.Count; $i++) {
  $pp = $ftr.Range.Paragraphs($i)
  Write-Output "Para $i len=$($pp.Range.End - $pp.Range.Start) text=[$($pp.Range.Text)]"
}
